$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# Model box (id=56)
$sh = Get-ShapeById $s 56
$sh.Left = 513.9732055664062
$sh.Width = 329.1040954589844

# Logic box (id=81)
$sh = Get-ShapeById $s 81
$sh.Left = 38.652679443359375
$sh.Top = 12.672284126281738
$sh.Width = 476.76513671875

# Rectangle 62 (:Address/BookParser) (id=16)
$sh = Get-ShapeById $s 16
$sh.Left = 264.0016784667969
$sh.Width = 107.99669647216797

# TextBox 78 (AddressBook) (id=79)
$sh = Get-ShapeById $s 79
$sh.Width = 121.1900863647461

# Rectangle 62 (VersionedAddressBook) (id=84)
$sh = Get-ShapeById $s 84
$sh.Left = 617.4044189453125
$sh.Width = 198.2838592529297

# Rectangle 85 (id=86)
$sh = Get-ShapeById $s 86
$sh.Left = 697.4873046875

# Rectangle 62 (id=40)
$sh = Get-ShapeById $s 40
$sh.Left = 542.8068237304688

# Straight Connector 45 (id=46)
$sh = Get-ShapeById $s 46
$sh.Left = 577.1947021484375

# Rectangle 48 (id=49)
$sh = Get-ShapeById $s 49
$sh.Left = 569.4022216796875

# Straight Arrow Connector 49 (id=50)
$sh = Get-ShapeById $s 50
$sh.Top = 232.9855194091797
$sh.Width = 122.32669830322266
$sh.Height = 1.2935433387756348
$sh.VerticalFlip = -1  # msoTrue

# TextBox 87 (ReadOnlyAddressBook) (id=88)
$sh = Get-ShapeById $s 88
$sh.Left = 717.5913696289062
$sh.Width = 188.40858459472656

# Straight Connector 88 (id=89)
$sh = Get-ShapeById $s 89
$sh.Left = 703.3011474609375

# TextBox 40 (id=41)
$sh = Get-ShapeById $s 41
$sh.Left = 630.8324584960938

# Curved Connector 12 (id=13)
$sh = Get-ShapeById $s 13
$sh.Left = 698.5401000976562

# Straight Arrow Connector 51 (id=52)
$sh = Get-ShapeById $s 52
$sh.Left = 581.7222290039062

# Straight Arrow Connector 59 (id=60)
$sh = Get-ShapeById $s 60
$sh.Left = 581.7222290039062

# Straight Arrow Connector 62 (id=63)
$sh = Get-ShapeById $s 63
$sh.Left = 447.10882568359375
$sh.Width = 128.94284057617188

# --- Text changes: AddressBook -> FinancialPlanner renaming ---

# Rectangle 62 ":Address" / "BookParser" -> ":Financial" / "PlannerParser" (id=16)
$sh = Get-ShapeById $s 16
$sh.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = ":Financial"
$sh.TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "PlannerParser"

# TextBox 78 "undo" + "AddressBook" + "()" -> "undo" + "FinancialPlanner" + "()" (id=79)
$sh = Get-ShapeById $s 79
$sh.TextFrame.TextRange.Paragraphs(1).Runs(2).Text = "FinancialPlanner"

# Rectangle 62 ":" + "VersionedAddressBook" -> ":" + "VersionedFinancialPlanner" (id=84)
$sh = Get-ShapeById $s 84
$sh.TextFrame.TextRange.Paragraphs(1).Runs(2).Text = "VersionedFinancialPlanner"

# TextBox 87 "resetData" + "(" + "ReadOnlyAddressBook" + ")" -> ... "ReadOnlyFinancialPlanner" ... (id=88)
$sh = Get-ShapeById $s 88
$sh.TextFrame.TextRange.Paragraphs(1).Runs(3).Text = "ReadOnlyFinancialPlanner"
